# Auto-generated COM-interop script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '66.782.22'
Set-TextValue 'E2' '  +2.07%  '
Set-TextValue 'D3' '3.272.83'
Set-TextValue 'E3' '  -0.51%  '
Set-TextValue 'D4' '0.998'
Set-TextValue 'E4' '  +0.08%  '
Set-TextValue 'D5' '569.53'
Set-TextValue 'E5' '  -1.43%  '
Set-TextValue 'D6' '175.66'
Set-TextValue 'E6' '  -3.88%  '
Set-TextValue 'E7' '  -0.08%  '
Set-TextValue 'D8' '0.582'
Set-TextValue 'E8' '  +2.21%  '
Set-TextValue 'D9' '3.267.42'
Set-TextValue 'E9' '  -0.48%  '
Set-TextValue 'D10' '0.173'
Set-TextValue 'E10' '  -0.83%  '
Set-TextValue 'D11' '0.570'
Set-TextValue 'E11' '  +0.03%  '
Set-TextValue 'D12' '45.47'
Set-TextValue 'E12' '  -2.00%  '
Set-TextValue 'D13' '0.0000268'
Set-TextValue 'E13' '  +2.04%  '
Set-TextValue 'D14' '692.26'
Set-TextValue 'E14' '  +9.55%  '
Set-TextValue 'D15' '3.792.14'
Set-TextValue 'E15' '  -0.56%  '
Set-TextValue 'D16' '8.29'
Set-TextValue 'E16' '  -1.13%  '
Set-TextValue 'D17' '66.820.41'
Set-TextValue 'E17' '  +1.90%  '
Set-TextValue 'E18' '  +1.08%  '
Set-TextValue 'D19' '3.274.24'
Set-TextValue 'E19' '  -0.43%  '
Set-TextValue 'D20' '17.29'
Set-TextValue 'E20' '  -1.89%  '
Set-TextValue 'D21' '10.73'
Set-TextValue 'E21' '  -1.62%  '
Set-TextValue 'D22' '0.885'
Set-TextValue 'E22' '  +0.00%  '
Set-TextValue 'D23' '16.90'
Set-TextValue 'E23' '  -5.72%  '
Set-TextValue 'D24' '5.13'
Set-TextValue 'E24' '  +3.44%  '
Set-TextValue 'D25' '97.62'
Set-TextValue 'E25' '  -2.51%  '
Set-TextValue 'D26' '3.86'
Set-TextValue 'E26' '  -2.12%  '
Set-TextValue 'D27' '2.70'
Set-TextValue 'E27' '  -1.02%  '
Set-TextValue 'D28' '9.28'
Set-TextValue 'E28' '  -0.58%  '
Set-TextValue 'D29' '32.79'
Set-TextValue 'E29' '  +7.11%  '
Set-TextValue 'D30' '8.40'
Set-TextValue 'E30' '  +0.88%  '
Set-TextValue 'D31' '6.76'
Set-TextValue 'E31' '  +4.50%  '
Set-TextValue 'D32' '574.75'
Set-TextValue 'E32' '  -0.02%  '
Set-TextValue 'D33' '3.869.36'
Set-TextValue 'E33' '  +0.69%  '
Set-TextValue 'D34' '10.77'
Set-TextValue 'E34' '  -0.58%  '
Set-TextValue 'E35' '  +0.13%  '
Set-TextValue 'E36' '  -0.03%  '
Set-TextValue 'D37' '55.32'
Set-TextValue 'E37' '  -0.24%  '
Set-TextValue 'D38' '3.29'
Set-TextValue 'E38' '  -10.32%  '
Set-TextValue 'E39' '  +2.44%  '
Set-TextValue 'D40' '2.60'
Set-TextValue 'E40' '  +0.57%  '
Set-TextValue 'E41' '  -1.76%  '
Set-TextValue 'D42' '31.65'
Set-TextValue 'E42' '  -2.14%  '
Set-TextValue 'D43' '0.0₃0668'
Set-TextValue 'E43' '  -1.44%  '
Set-TextValue 'D44' '3.03'
Set-TextValue 'E44' '  -2.36%  '
Set-TextValue 'D45' '0.327'
Set-TextValue 'E45' '  -1.17%  '
Set-TextValue 'D46' '0.0405'
Set-TextValue 'E46' '  +0.40%  '
Set-TextValue 'E47' '  +0.66%  '
Set-TextValue 'E48' '  +0.22%  '
Set-TextValue 'D51' '129.15'
Set-TextValue 'E51' '  -0.15%  '

# Rows 49/50: Mantle and ThetaToken swap ranking positions with updated figures
Set-TextValue 'B49' 'ThetaToken'
Set-TextValue 'C49' 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D49' '2.52'
Set-TextValue 'E49' '  +0.90%  '

Set-TextValue 'B50' 'Mantle'
Set-TextValue 'C50' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D50' '1.36'
Set-TextValue 'E50' '  +7.97%  '
